$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for a new "Periodo Mora" row: insert a blank row at 21, then
#    paste down the formatting (borders/number formats) of row 20 (the last
#    of the "middle style" data rows) so the new row matches the table style.
#    This pushes the previous row 21 (and everything below it, including the
#    signature block) down by one row.
$ws.Rows("21:21").Insert() | Out-Null
$ws.Rows("20:20").Copy() | Out-Null
$ws.Rows("21:21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2) Fill in the new row with the same worker/values as the rest of the
#    table, just a different period.
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "11203043"
$ws.Range("D21").Value = "SIMONIDES VASCO VELEZ"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 360000
$ws.Range("G21").Value = 9000000

# 3) Refresh the "Periodo Mora" (period) column for the rest of the worker
#    table: the periods now run in ascending order 2502..2508, with row 22
#    (the old closing row, now pushed down) becoming the brand-new period
#    2508 entry.
$ws.Range("E16").Value = "2502"
$ws.Range("E17").Value = "2503"
$ws.Range("E18").Value = "2504"
$ws.Range("E19").Value = "2505"
$ws.Range("E20").Value = "2506"
$ws.Range("E22").Value = "2508"

# 4) Update the summary totals: total overdue amount and period count.
$ws.Range("E11").Value = 2520000
$ws.Range("F13").Value = 7
